$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header column U1
$ws.Range("U1").Value = "Fertilizer Recommendation"

# Updated row 2 values
$ws.Range("B2").Value = "24-04-2024"
$ws.Range("C2").Value = 415
$ws.Range("D2").Value = 321
$ws.Range("E2").Value = "akakjdklasd"
$ws.Range("F2").Value = 3
$ws.Range("H2").Value = 25
$ws.Range("I2").Value = "654asd65asd"

# J2 (Mobile No.) is a purely-numeric-looking string; force text format
# so it is stored as text rather than auto-converted to a number.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "8456231231"
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 250
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 80
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 25
$ws.Range("Q2").Value = 30
$ws.Range("R2").Value = 20
$ws.Range("S2").Value = 0.3625365324113583
$ws.Range("T2").Value = "Grow Millets (Sorghum, Pearl millet), Pulses (Pigeon pea, Chickpea), and Oilseeds (Safflower, Castor)."
$ws.Range("U2").Value = "Apply organic amendments like Compost (3-5 tonnes/ha), Vermicompost (1.5-2.5 tonnes/ha), or Well-decomposed Farmyard manure (7.5-10 tonnes/ha). Use biofertilizers like Rhizobium (200-300 g/ha), Azotobacter (200-300 g/ha), and Phosphate Solubilizing Bacteria (PSB) (500-750 g/ha). Apply chemical fertilizers at 50% of the recommended dose based on soil test results."
